$wb = $excel.ActiveWorkbook

# --- Sheet "2o Parcial", row 8 (García Sánchez Magda Bexabe / 2ALCV) ---
$ws2 = $wb.Worksheets.Item("2o Parcial")
$ws2.Range("E8").Value = 22
$ws2.Range("F8").Value = 10
$ws2.Range("G8").Value = 68.75
$ws2.Range("H8").Value = 31.25
$ws2.Range("I8").Value = 8.1
$ws2.Range("J8").Value = 10
$ws2.Range("K8").Value = 31.25

# --- Sheet "3er Parcial", row 8 (same teacher/group) ---
$ws3 = $wb.Worksheets.Item("3er Parcial")
$ws3.Range("I8").Value = 8.1
